# Scheduled runner update: refresh market-board price snapshots (and the
# dependent Leve profit totals) for several crafting-leve rows across sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1124.9474
$ws.Range("J129").Value = 1810
$ws.Range("L129").Value = 5430
$ws.Range("N129").Value = -15430
$ws.Range("H137").Value = 2434.913
$ws.Range("I137").Value = 2042.8572
$ws.Range("J137").Value = 2606.4375
$ws.Range("K137").Value = 6128.571599999999
$ws.Range("L137").Value = 7819.3125
$ws.Range("M137").Value = -3578.571599999999
$ws.Range("N137").Value = -12919.3125
$ws.Range("H138").Value = 2797.9033
$ws.Range("I138").Value = 7166.6665
$ws.Range("J138").Value = 2575.7627
$ws.Range("K138").Value = 21499.9995
$ws.Range("L138").Value = 7727.288100000001
$ws.Range("M138").Value = -16359.9995
$ws.Range("N138").Value = -18007.2881

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 30941.715
$ws.Range("J7").Value = 30941.715
$ws.Range("L7").Value = 30941.715
$ws.Range("N7").Value = -31169.715
$ws.Range("H22").Value = 7293.5
$ws.Range("I22").Value = 7750.4
$ws.Range("J22").Value = 5009
$ws.Range("K22").Value = 7750.4
$ws.Range("L22").Value = 5009
$ws.Range("M22").Value = -7451.4
$ws.Range("N22").Value = -5607
$ws.Range("H24").Value = 28195.643
$ws.Range("J24").Value = 28195.643
$ws.Range("L24").Value = 28195.643
$ws.Range("N24").Value = -28943.643
$ws.Range("H96").Value = 25297
$ws.Range("J96").Value = 25297
$ws.Range("L96").Value = 25297
$ws.Range("N96").Value = -30789
$ws.Range("H100").Value = 28195.643
$ws.Range("J100").Value = 28195.643
$ws.Range("L100").Value = 28195.643
$ws.Range("N100").Value = -30359.643
$ws.Range("H109").Value = 8806.076999999999
$ws.Range("J109").Value = 8806.076999999999
$ws.Range("L109").Value = 8806.076999999999
$ws.Range("N109").Value = -11580.077
$ws.Range("H112").Value = 16395.666
$ws.Range("J112").Value = 16395.666
$ws.Range("L112").Value = 16395.666
$ws.Range("N112").Value = -19349.666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H108").Value = 31149.2
$ws.Range("J108").Value = 31149.2
$ws.Range("L108").Value = 31149.2
$ws.Range("N108").Value = -38829.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 177.8
$ws.Range("I7").Value = 40
$ws.Range("J7").Value = 269.66666
$ws.Range("K7").Value = 40
$ws.Range("L7").Value = 269.66666
$ws.Range("M7").Value = 73
$ws.Range("N7").Value = -495.66666
$ws.Range("H31").Value = 2548.4524
$ws.Range("I31").Value = 1540
$ws.Range("K31").Value = 1540
$ws.Range("M31").Value = -1245
$ws.Range("H34").Value = 2548.4524
$ws.Range("I34").Value = 1540
$ws.Range("K34").Value = 1540
$ws.Range("M34").Value = -1338
$ws.Range("H53").Value = 44748.75
$ws.Range("J53").Value = 44748.75
$ws.Range("L53").Value = 44748.75
$ws.Range("N53").Value = -45962.75
$ws.Range("H99").Value = 1894.1765
$ws.Range("I99").Value = 1690.5
$ws.Range("J99").Value = 2075.2222
$ws.Range("K99").Value = 1690.5
$ws.Range("L99").Value = 2075.2222
$ws.Range("M99").Value = -192.5
$ws.Range("N99").Value = -5071.2222
$ws.Range("H106").Value = 27000
$ws.Range("J106").Value = 27000
$ws.Range("L106").Value = 27000
$ws.Range("N106").Value = -29524
$ws.Range("H111").Value = 32700
$ws.Range("J111").Value = 32700
$ws.Range("L111").Value = 32700
$ws.Range("N111").Value = -40880
$ws.Range("H122").Value = 2966.05
$ws.Range("I122").Value = 1512.6666
$ws.Range("J122").Value = 4155.1816
$ws.Range("K122").Value = 4537.9998
$ws.Range("L122").Value = 12465.5448
$ws.Range("M122").Value = -2087.9998
$ws.Range("N122").Value = -17365.5448
$ws.Range("H126").Value = 1894.1765
$ws.Range("I126").Value = 1690.5
$ws.Range("J126").Value = 2075.2222
$ws.Range("K126").Value = 5071.5
$ws.Range("L126").Value = 6225.6666
$ws.Range("M126").Value = -2601.5
$ws.Range("N126").Value = -11165.6666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 662.4286
$ws.Range("I113").Value = 606.6667
$ws.Range("J113").Value = 704.25
$ws.Range("K113").Value = 1820.0001
$ws.Range("L113").Value = 2112.75
$ws.Range("M113").Value = 349.9999
$ws.Range("N113").Value = -6452.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 53164.43
$ws.Range("I70").Value = 106450.8
$ws.Range("J70").Value = 4722.273
$ws.Range("K70").Value = 106450.8
$ws.Range("L70").Value = 4722.273
$ws.Range("M70").Value = -106180.8
$ws.Range("N70").Value = -5262.273
$ws.Range("H73").Value = 53164.43
$ws.Range("I73").Value = 106450.8
$ws.Range("J73").Value = 4722.273
$ws.Range("K73").Value = 106450.8
$ws.Range("L73").Value = 4722.273
$ws.Range("M73").Value = -105514.8
$ws.Range("N73").Value = -6594.273
$ws.Range("H95").Value = 8577.125
$ws.Range("J95").Value = 8577.125
$ws.Range("L95").Value = 8577.125
$ws.Range("N95").Value = -14069.125
$ws.Range("H103").Value = 28500
$ws.Range("J103").Value = 28500
$ws.Range("L103").Value = 28500
$ws.Range("N103").Value = -30844
$ws.Range("H122").Value = 2349.5833
$ws.Range("I122").Value = 2586.7144
$ws.Range("J122").Value = 2017.6
$ws.Range("K122").Value = 7760.1432
$ws.Range("L122").Value = 6052.799999999999
$ws.Range("M122").Value = -5310.1432
$ws.Range("N122").Value = -10952.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1108.8889
$ws.Range("J46").Value = 1260
$ws.Range("L46").Value = 1260
$ws.Range("N46").Value = -1636
$ws.Range("H93").Value = 2174.2666
$ws.Range("I93").Value = 2023.3334
$ws.Range("J93").Value = 2400.6667
$ws.Range("K93").Value = 2023.3334
$ws.Range("L93").Value = 2400.6667
$ws.Range("M93").Value = -775.3334
$ws.Range("N93").Value = -4896.6667
$ws.Range("H94").Value = 17220
$ws.Range("J94").Value = 17220
$ws.Range("L94").Value = 17220
$ws.Range("N94").Value = -18572
$ws.Range("H100").Value = 2442.8572
$ws.Range("I100").Value = 1600
$ws.Range("J100").Value = 3566.6667
$ws.Range("K100").Value = 1600
$ws.Range("L100").Value = 3566.6667
$ws.Range("M100").Value = -1059
$ws.Range("N100").Value = -4648.6667
$ws.Range("H132").Value = 4851.619
$ws.Range("I132").Value = 5030.375
$ws.Range("J132").Value = 4279.6
$ws.Range("K132").Value = 15091.125
$ws.Range("L132").Value = 12838.8
$ws.Range("M132").Value = -12561.125
$ws.Range("N132").Value = -17898.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H116").Value = 30199.666
$ws.Range("J116").Value = 30199.666
$ws.Range("L116").Value = 30199.666
$ws.Range("N116").Value = -39377.666
$ws.Range("H122").Value = 2300.5
$ws.Range("I122").Value = 2874.1333
$ws.Range("J122").Value = 1638.6154
$ws.Range("K122").Value = 8622.3999
$ws.Range("L122").Value = 4915.8462
$ws.Range("M122").Value = -6172.3999
$ws.Range("N122").Value = -9815.8462
